# Revert capacity chart to show kilowatts (not watts) on the y-axis.
#  - Data cells (and their shared number format) switch from whole-number
#    "#,##0" to one-decimal "#,##0.0" display.
#  - The "Solar" column's 2018/2019/2022/2023/2024 figures were stored in
#    watts (e.g. 7600) and are corrected to kilowatts (7.6).
#  - The value axis title/number format on the chart go back to a plain
#    kilowatts label instead of the "K"-suffixed watts format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data: convert the Solar column's watt figures to kilowatts ---
$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 6
$ws.Range("E24").Value = 7.6
$ws.Range("E25").Value = 11.6
$ws.Range("E26").Value = 3.48

# Show one decimal place across the data grid (B2:G26), matching the new
# kilowatt-scale values.
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- Chart: fix the value axis title + number format, and refresh the ---
# --- cached series values so the chart matches the worksheet data.    ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$valueAxis = $chart.Axes(2)
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"

$series = $chart.SeriesCollection(4)
$series.Values = $ws.Range("E2:E26")
